$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 1.77

# Row 3
$ws.Range("V3").Value = 1.63

# Row 4
$ws.Range("V4").Value = 1.63

# Row 5
$ws.Range("G5").Value = 2.8
$ws.Range("I5").Value = 2.4
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.25
$ws.Range("X5").Value = 15
$ws.Range("AN5").Value = 5
$ws.Range("BA5").Value = 51

# Row 9
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
